$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Fix site: ... About page: ... Add short bio" block.
#    That block sits between the empty centered paragraph right after
#    the title and the "Host the site:" paragraph.
# ------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.TrimEnd([char]13, [char]7) -eq "Fix site:") {
        $startPara = $i
    }
    if ($txt.TrimEnd([char]13, [char]7) -eq "Add short bio") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End
    $killRange = $d.Range($rangeStart, $rangeEnd)
    $killRange.Delete()
}

# ------------------------------------------------------------------
# 2) Split the run " PySpark, Airflow etc ... necessary" into three
#    runs and wrap "PySpark" with spell-check proof-error markers,
#    matching what Word produces after nudging the text next to it.
# ------------------------------------------------------------------
$dash = [char]0x2013
$needle = " PySpark, Airflow etc " + $dash + " probably not 100% necessary"

$findRange = $d.Content
$found = $findRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($findRange.Start -ge $candidate.Range.Start -and $findRange.Start -lt $candidate.Range.End) {
            $para = $candidate
        }
    }
    $rsidTarget = "00023C84"
    $rsidPr = "00C61810"

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="52EC6754" w14:textId="34BA144A" w:rsidR="001B00D4" w:rsidRPr="' + $rsidPr + '" w:rsidRDefault="001B00D4" w:rsidP="000A1B3A">' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>' + `
        '<w:r w:rsidRPr="' + $rsidPr + '"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Practice project using deployment?</w:t></w:r>' + `
        '<w:r w:rsidR="' + $rsidTarget + '" w:rsidRPr="' + $rsidPr + '"><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r w:rsidR="' + $rsidTarget + '" w:rsidRPr="' + $rsidPr + '"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>PySpark</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r w:rsidR="' + $rsidTarget + '" w:rsidRPr="' + $rsidPr + '"><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">, Airflow etc ' + $dash + ' probably not 100% necessary</w:t></w:r>' + `
        '<w:r w:rsidR="00C61810" w:rsidRPr="' + $rsidPr + '"><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> at this stage</w:t></w:r>' + `
        '</w:p>'

    if ($para -ne $null) {
        $para.Range.InsertXML($xml)
    }
}
